$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "Nombre" -> "Nombre de producto"
Replace-Text "Nombre" "Nombre de producto"

# 2. " del producto: Místico Spice Premium Chai Tea" -> ": Té chai premium Mystic Spice"
Replace-Text " del producto: Místico Spice Premium Chai Tea" ": Té chai premium Mystic Spice"

# 3. Bold "Principales características:" run (b val=0 -> b)
# 4. "Principales características:" -> "Características principales:"
Replace-Text "Principales características:" "Características principales:"

# 5. "Mezcla" -> "Mezcla auténtica"
Replace-Text "Mezcla" "Mezcla auténtica"

# 6.
Replace-Text " auténtica: Nuestra chai es una mezcla armónica de hojas de té negro premium y una selección de especias molidas, incluyendo canela, cardamomo, cloves, jengibre y pimienta negra." ": nuestro chai es una mezcla armoniosa de hojas de té negro premium y una selección de especias molidas, incluyendo canela, cardamomo, clavo, jengibre y pimienta negra."

# 7. "Ingredientes" -> "Ingredientes beneficiosos para la salud"
Replace-Text "Ingredientes" "Ingredientes beneficiosos para la salud"

# 8.
Replace-Text " de mejora de la salud: Cada ingrediente del Místico Spice Chai Tea se elige para sus beneficios naturales para la salud." ": cada ingrediente del té chai premium Mystic Spice se elige por sus beneficios naturales para la salud."

# 9.
Replace-Text "El jengibre y el cardamomo ayudan a la digestión, la canela ayuda a regular el azúcar en sangre y los clavos agregan un impulso de antioxidantes." "El jengibre y el cardamomo ayudan a la digestión, la canela ayuda a regular el azúcar en sangre y el clavo aumento los antioxidantes."

# 10. "Rico Aroma y Sabor" -> "Aroma y sabor intensos"
Replace-Text "Rico Aroma y Sabor" "Aroma y sabor intensos"

# 11.
Replace-Text ": El aroma cálido, picante y profundo, vigorizante sabor de nuestra chai hacen que sea la bebida perfecta para comenzar su día o relajarse por la noche." ": el aroma cálido y especiado, y el sabor profundo y vigorizante de nuestro chai hacen que sea la bebida perfecta para comenzar el día o relajarse por la noche."

# 12. "Opciones" -> "Opciones versátiles de preparación"
Replace-Text "Opciones" "Opciones versátiles de preparación"

# 13.
Replace-Text " versátiles de preparación: Ya sea que amas tu chai vaporing caliente, como un refrescante té helado, o como una latte cremosa, nuestra mezcla es lo suficientemente versátil como para adaptarte a cualquier preferencia." ": ya sea que te guste tu chai bien caliente, o prefieras un refrescante té helado, o un latte cremoso, nuestra mezcla es lo suficientemente versátil como para adaptarse a cualquier preferencia."

# 14. "Origen" -> "Origen sostenible"
Replace-Text "Origen" "Origen sostenible"

# 15.
Replace-Text " sostenible: Comprometidos con la sostenibilidad, originamos nuestros ingredientes de granjas a pequeña escala que practican la agricultura ecológica, garantizando no sólo la mejor calidad, sino también el bienestar de nuestro planeta." ": al estar comprometidos con la sostenibilidad, obtenemos nuestros ingredientes de pequeñas explotaciones que practican la agricultura ecológica, garantizando no solo la mejor calidad, sino también el bienestar de nuestro planeta."

# 16. "Empaquetado" -> "Envase elegante"
Replace-Text "Empaquetado" "Envase elegante"

# 17.
Replace-Text " elegante: El té de Spice Chai místico viene en un empaquetado elegante, ecológico, lo que lo convierte en un regalo ideal para los amantes del té o un lujoso trato para usted mismo." ": el té chai Mystic Spice viene en un envase elegante, ecológico, lo que lo convierte en el regalo ideal para los amantes del té o un capricho lujoso para ti mismo."

# 18. "Garantía" -> "Garantía de satisfacción del cliente"
Replace-Text "Garantía" "Garantía de satisfacción del cliente"

# 19.
Replace-Text " de satisfacción del cliente: Estamos detrás de nuestro producto y ofrecemos una garantía de satisfacción." ": respaldamos nuestro producto y ofrecemos una garantía de satisfacción."

# 20.
Replace-Text ": entusiastas del té, individuos conscientes de la salud, amantes de bebidas calientes, especiadas, y cualquier persona que busca explorar los ricos sabores de la chai india tradicional." ": los apasionados del té, las personas conscientes de la salud, los amantes de las bebidas calientes especiadas, y cualquier persona que desea explorar los sabores intensos del chai indio tradicional."

